{"js": "// Replace the date line and the 25 division problems/answers in the table\n// with their updated values, as described by the commit diff.\nconst replacements = [\n  [\"2023-11-29 Wednesday\", \"2023-11-30 Thursday\"],\n  [\"77\u00f74=19, 1\", \"35\u00f77=5, 0\"],\n  [\"78\u00f76=13, 0\", \"97\u00f78=12, 1\"],\n  [\"29\u00f79=3, 2\", \"11\u00f76=1, 5\"],\n  [\"64\u00f76=10, 4\", \"75\u00f79=8, 3\"],\n  [\"50\u00f79=5, 5\", \"34\u00f75=6, 4\"],\n  [\"21\u00f74=5, 1\", \"59\u00f74=14, 3\"],\n  [\"82\u00f76=13, 4\", \"47\u00f76=7, 5\"],\n  [\"79\u00f73=26, 1\", \"78\u00f73=26, 0\"],\n  [\"36\u00f79=4, 0\", \"22\u00f72=11, 0\"],\n  [\"67\u00f74=16, 3\", \"84\u00f78=10, 4\"],\n  [\"29\u00f74=7, 1\", \"51\u00f72=25, 1\"],\n  [\"98\u00f77=14, 0\", \"90\u00f78=11, 2\"],\n  [\"60\u00f72=30, 0\", \"69\u00f72=34, 1\"],\n  [\"44\u00f78=5, 4\", \"35\u00f73=11, 2\"],\n  [\"99\u00f78=12, 3\", \"59\u00f72=29, 1\"],\n  [\"46\u00f79=5, 1\", \"13\u00f77=1, 6\"],\n  [\"46\u00f78=5, 6\", \"71\u00f78=8, 7\"],\n  [\"82\u00f79=9, 1\", \"67\u00f73=22, 1\"],\n  [\"15\u00f76=2, 3\", \"17\u00f79=1, 8\"],\n  [\"99\u00f72=49, 1\", \"63\u00f72=31, 1\"],\n  [\"40\u00f76=6, 4\", \"32\u00f72=16, 0\"],\n  [\"31\u00f75=6, 1\", \"99\u00f78=12, 3\"],\n  [\"83\u00f79=9, 2\", \"92\u00f76=15, 2\"],\n  [\"38\u00f79=4, 2\", \"51\u00f76=8, 3\"],\n  [\"93\u00f73=31, 0\", \"18\u00f77=2, 4\"],\n];\n\n// Each old value occurs exactly once in the document, so searching and\n// replacing one-at-a-time (in document order) is safe even though some\n// new values coincide with old values used elsewhere (e.g. \"99\u00f78=12, 3\"\n// is both an old value near the top and a new value further down).\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 division problems/answers in the table\n# with their updated values, as described by the commit diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2023-11-29 Wednesday\", \"2023-11-30 Thursday\"),\n  @(\"77\u00f74=19, 1\", \"35\u00f77=5, 0\"),\n  @(\"78\u00f76=13, 0\", \"97\u00f78=12, 1\"),\n  @(\"29\u00f79=3, 2\", \"11\u00f76=1, 5\"),\n  @(\"64\u00f76=10, 4\", \"75\u00f79=8, 3\"),\n  @(\"50\u00f79=5, 5\", \"34\u00f75=6, 4\"),\n  @(\"21\u00f74=5, 1\", \"59\u00f74=14, 3\"),\n  @(\"82\u00f76=13, 4\", \"47\u00f76=7, 5\"),\n  @(\"79\u00f73=26, 1\", \"78\u00f73=26, 0\"),\n  @(\"36\u00f79=4, 0\", \"22\u00f72=11, 0\"),\n  @(\"67\u00f74=16, 3\", \"84\u00f78=10, 4\"),\n  @(\"29\u00f74=7, 1\", \"51\u00f72=25, 1\"),\n  @(\"98\u00f77=14, 0\", \"90\u00f78=11, 2\"),\n  @(\"60\u00f72=30, 0\", \"69\u00f72=34, 1\"),\n  @(\"44\u00f78=5, 4\", \"35\u00f73=11, 2\"),\n  @(\"99\u00f78=12, 3\", \"59\u00f72=29, 1\"),\n  @(\"46\u00f79=5, 1\", \"13\u00f77=1, 6\"),\n  @(\"46\u00f78=5, 6\", \"71\u00f78=8, 7\"),\n  @(\"82\u00f79=9, 1\", \"67\u00f73=22, 1\"),\n  @(\"15\u00f76=2, 3\", \"17\u00f79=1, 8\"),\n  @(\"99\u00f72=49, 1\", \"63\u00f72=31, 1\"),\n  @(\"40\u00f76=6, 4\", \"32\u00f72=16, 0\"),\n  @(\"31\u00f75=6, 1\", \"99\u00f78=12, 3\"),\n  @(\"83\u00f79=9, 2\", \"92\u00f76=15, 2\"),\n  @(\"38\u00f79=4, 2\", \"51\u00f76=8, 3\"),\n  @(\"93\u00f73=31, 0\", \"18\u00f77=2, 4\")\n)\n\n# Each old value occurs exactly once in the document, so a single\n# wdReplaceOne Find/Replace per pair (walked in document order) is safe\n# even though some new values coincide with old values used elsewhere\n# (e.g. \"99\u00f78=12, 3\" is both an old value near the top and a new value\n# further down the table).\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $range = $d.Content\n  $found = $range.Find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 1)\n  if (-not $found) {\n    Write-Output \"NOT FOUND: $old\"\n  }\n}\n"}
